# Storage component diagram: fix component naming
#
# The "PresentationSectionRepository" rounded-rectangle box was too
# narrow for its (renamed/corrected) label, so it is widened; the
# dashed connector that feeds into it from the right is shortened and
# re-anchored to meet the box's new (shifted) connection point, and its
# elbow bend ("adj1") is recomputed to keep the bend midway along the
# shorter run.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "Rounded Rectangle 18" - the PresentationSectionRepository box.
# Only its width changes (2590800 EMU -> 2971800 EMU); position/height
# are untouched.
$repositoryBox = $s.Shapes.Item(14)
$repositoryBox.Width = 234

# "Straight Arrow Connector 161" that lands on the repository box
# (stCxn id=19 idx=3 -> endCxn id=41 idx=0). Its left edge moves right
# and its width shrinks to match, while top/height stay the same; the
# bentConnector3 "adj1" guide (the fraction along the bend) is updated
# from 84250/100000 to 66104/100000.
$connector = $s.Shapes.Item(25)
$connector.Left = 544.1874803149607
$connector.Width = 64.4457513648294
$connector.Adjustments.Item(1) = 0.66104
